$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 297.2
$ws.Range("I6").Value = 290
$ws.Range("K6").Value = 870
$ws.Range("M6").Value = -758
$ws.Range("H53").Value = 1134.7273
$ws.Range("I53").Value = 1483.2858
$ws.Range("J53").Value = 524.75
$ws.Range("K53").Value = 1483.2858
$ws.Range("L53").Value = 524.75
$ws.Range("M53").Value = -846.2858000000001
$ws.Range("N53").Value = -1798.75
$ws.Range("H62").Value = 5289.625
$ws.Range("I62").Value = 4447.3335
$ws.Range("K62").Value = 4447.3335
$ws.Range("M62").Value = -3823.3335
$ws.Range("H65").Value = 5289.625
$ws.Range("I65").Value = 4447.3335
$ws.Range("K65").Value = 22236.6675
$ws.Range("M65").Value = -19116.6675
$ws.Range("H107").Value = 56820308
$ws.Range("I107").Value = 12502338
$ws.Range("J107").Value = 500000000
$ws.Range("K107").Value = 12502338
$ws.Range("L107").Value = 500000000
$ws.Range("M107").Value = -12500418
$ws.Range("N107").Value = -500003840
$ws.Range("H117").Value = 14900
$ws.Range("J117").Value = 14900
$ws.Range("L117").Value = 14900
$ws.Range("N117").Value = -24078
$ws.Range("H128").Value = 95000
$ws.Range("J128").Value = 95000
$ws.Range("L128").Value = 95000
$ws.Range("N128").Value = -104960
$ws.Range("H135").Value = 455012
$ws.Range("I135").Value = 476622.66
$ws.Range("K135").Value = 4289603.939999999
$ws.Range("M135").Value = -4287068.939999999
$ws.Range("H138").Value = 1003484.5
$ws.Range("I138").Value = 1869.6666
$ws.Range("J138").Value = 1643861.2
$ws.Range("K138").Value = 5608.9998
$ws.Range("L138").Value = 4931583.6
$ws.Range("M138").Value = -468.9997999999996
$ws.Range("N138").Value = -4941863.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 133
$ws.Range("I4").Value = 99.5
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 99.5
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 16.5
$ws.Range("N4").Value = -432
$ws.Range("H5").Value = 14617
$ws.Range("I5").Value = 16986.5
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 16986.5
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -16874.5
$ws.Range("N5").Value = -624
$ws.Range("H32").Value = 5619.9214
$ws.Range("I32").Value = 5119.706
$ws.Range("K32").Value = 5119.706
$ws.Range("M32").Value = -4832.706
$ws.Range("H74").Value = 20068.875
$ws.Range("I74").Value = 30873.854
$ws.Range("K74").Value = 30873.854
$ws.Range("M74").Value = -29999.854
$ws.Range("H77").Value = 20068.875
$ws.Range("I77").Value = 30873.854
$ws.Range("K77").Value = 154369.27
$ws.Range("M77").Value = -150001.27
$ws.Range("H122").Value = 2150.7727
$ws.Range("I122").Value = 1792.7778
$ws.Range("K122").Value = 5378.3334
$ws.Range("M122").Value = -2928.3334
$ws.Range("H132").Value = 3961.0698
$ws.Range("I132").Value = 2031.56
$ws.Range("K132").Value = 6094.68
$ws.Range("M132").Value = -3564.68
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 14617
$ws.Range("I4").Value = 16986.5
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 16986.5
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -16871.5
$ws.Range("N4").Value = -630
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4392.927
$ws.Range("I31").Value = 1752.0435
$ws.Range("J31").Value = 7767.3887
$ws.Range("K31").Value = 1752.0435
$ws.Range("L31").Value = 7767.3887
$ws.Range("M31").Value = -1457.0435
$ws.Range("N31").Value = -8357.3887
$ws.Range("H34").Value = 4392.927
$ws.Range("I34").Value = 1752.0435
$ws.Range("J34").Value = 7767.3887
$ws.Range("K34").Value = 1752.0435
$ws.Range("L34").Value = 7767.3887
$ws.Range("M34").Value = -1550.0435
$ws.Range("N34").Value = -8171.3887
$ws.Range("H50").Value = 24428.572
$ws.Range("I50").Value = 5250
$ws.Range("J50").Value = 50000
$ws.Range("K50").Value = 5250
$ws.Range("L50").Value = 50000
$ws.Range("M50").Value = -4625
$ws.Range("N50").Value = -51250
$ws.Range("H51").Value = 45291.332
$ws.Range("J51").Value = 46487
$ws.Range("L51").Value = 46487
$ws.Range("N51").Value = -47959
$ws.Range("H61").Value = 45291.332
$ws.Range("J61").Value = 46487
$ws.Range("L61").Value = 46487
$ws.Range("N61").Value = -47183
$ws.Range("H132").Value = 3214.5156
$ws.Range("I132").Value = 1840.3541
$ws.Range("K132").Value = 5521.0623
$ws.Range("M132").Value = -2991.0623
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1240.4348
$ws.Range("I5").Value = 972.5294
$ws.Range("J5").Value = 1999.5
$ws.Range("K5").Value = 2917.5882
$ws.Range("L5").Value = 5998.5
$ws.Range("M5").Value = -2805.5882
$ws.Range("N5").Value = -6222.5
$ws.Range("H23").Value = 262.7
$ws.Range("J23").Value = 411.6
$ws.Range("L23").Value = 1234.8
$ws.Range("N23").Value = -1704.8
$ws.Range("H68").Value = 25004516
$ws.Range("J68").Value = 6694.4
$ws.Range("L68").Value = 20083.2
$ws.Range("N68").Value = -21705.2
$ws.Range("H71").Value = 25004516
$ws.Range("J71").Value = 6694.4
$ws.Range("L71").Value = 60249.6
$ws.Range("N71").Value = -68361.60000000001
$ws.Range("H107").Value = 13751096
$ws.Range("J107").Value = 17501340
$ws.Range("L107").Value = 52504020
$ws.Range("N107").Value = -52507860
$ws.Range("H113").Value = 4484.143
$ws.Range("J113").Value = 7100
$ws.Range("L113").Value = 21300
$ws.Range("N113").Value = -25640
$ws.Range("H122").Value = 1415982.5
$ws.Range("I122").Value = 2830062.5
$ws.Range("J122").Value = 1902.5
$ws.Range("K122").Value = 25470562.5
$ws.Range("L122").Value = 17122.5
$ws.Range("M122").Value = -25468112.5
$ws.Range("N122").Value = -22022.5
$ws.Range("H135").Value = 1240.4348
$ws.Range("I135").Value = 972.5294
$ws.Range("J135").Value = 1999.5
$ws.Range("K135").Value = 8752.7646
$ws.Range("L135").Value = 17995.5
$ws.Range("M135").Value = -6217.7646
$ws.Range("N135").Value = -23065.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = ""
$ws.Range("H47").Value = 17900
$ws.Range("J47").Value = 17900
$ws.Range("L47").Value = 17900
$ws.Range("N47").Value = -19036
$ws.Range("H102").Value = 1420.8889
$ws.Range("I102").Value = 1420.8889
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1420.8889
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 201.1111000000001
$ws.Range("N102").Value = ""
$ws.Range("H113").Value = 3159.0833
$ws.Range("I113").Value = 2705.2856
$ws.Range("J113").Value = 3794.4
$ws.Range("K113").Value = 2705.2856
$ws.Range("L113").Value = 3794.4
$ws.Range("M113").Value = -535.2856000000002
$ws.Range("N113").Value = -8134.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7858.9546
$ws.Range("I7").Value = 4857.6924
$ws.Range("J7").Value = 12194.111
$ws.Range("K7").Value = 4857.6924
$ws.Range("L7").Value = 12194.111
$ws.Range("M7").Value = -4745.6924
$ws.Range("N7").Value = -12418.111
$ws.Range("H16").Value = 1181.6154
$ws.Range("I16").Value = 1181.6154
$ws.Range("K16").Value = 1181.6154
$ws.Range("M16").Value = -1011.6154
$ws.Range("H122").Value = 4191.8887
$ws.Range("I122").Value = 3095.818
$ws.Range("K122").Value = 9287.454000000002
$ws.Range("M122").Value = -6837.454000000002
$ws.Range("H126").Value = 7858.9546
$ws.Range("I126").Value = 4857.6924
$ws.Range("J126").Value = 12194.111
$ws.Range("K126").Value = 14573.0772
$ws.Range("L126").Value = 36582.333
$ws.Range("M126").Value = -12103.0772
$ws.Range("N126").Value = -41522.333
$ws.Range("H132").Value = 17247284
$ws.Range("I132").Value = 29415648
$ws.Range("K132").Value = 88246944
$ws.Range("M132").Value = -88244414
$ws.Range("H136").Value = 9036.691999999999
$ws.Range("I136").Value = 3555.1482
$ws.Range("K136").Value = 10665.4446
$ws.Range("M136").Value = -8115.444600000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 43192.785
$ws.Range("I62").Value = 53339.273
$ws.Range("J62").Value = 5989
$ws.Range("K62").Value = 53339.273
$ws.Range("L62").Value = 5989
$ws.Range("M62").Value = -52715.273
$ws.Range("N62").Value = -7237
$ws.Range("H65").Value = 43192.785
$ws.Range("I65").Value = 53339.273
$ws.Range("J65").Value = 5989
$ws.Range("K65").Value = 266696.365
$ws.Range("L65").Value = 29945
$ws.Range("M65").Value = -263576.365
$ws.Range("N65").Value = -36185
$ws.Range("H96").Value = 2084.6667
$ws.Range("J96").Value = 1951.6
$ws.Range("L96").Value = 1951.6
$ws.Range("N96").Value = -4697.6
$ws.Range("H107").Value = 9804744
$ws.Range("I107").Value = 526.86664
$ws.Range("J107").Value = 17544914
$ws.Range("K107").Value = 1580.59992
$ws.Range("L107").Value = 52634742
$ws.Range("M107").Value = 339.4000800000001
$ws.Range("N107").Value = -52638582
$ws.Range("H132").Value = 6320.606
$ws.Range("I132").Value = 5824.7915
$ws.Range("K132").Value = 17474.3745
$ws.Range("M132").Value = -14944.3745
$ws.Range("H136").Value = 13119729
$ws.Range("I136").Value = 17858720
$ws.Range("J136").Value = 482421
$ws.Range("K136").Value = 53576160
$ws.Range("L136").Value = 1447263
$ws.Range("M136").Value = -53573610
$ws.Range("N136").Value = -1452363
